$d = $word.ActiveDocument

function Replace-One($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Note: "481x5=" occurs once in the original document; another cell is also being
# changed TO "481x5=" (272x8= -> 481x5=) later. To avoid the newly written "481x5="
# being matched again, update the pre-existing "481x5=" occurrence first.
Replace-One "481×5=" "374×6="

Replace-One "635×7=" "890×6="
Replace-One "258×9=" "779×4="
Replace-One "289×5=" "555×4="
Replace-One "272×8=" "481×5="
Replace-One "826×4=" "209×5="
Replace-One "116×4=" "324×9="
Replace-One "267×6=" "731×5="
Replace-One "419×4=" "377×7="
Replace-One "900×5=" "845×4="
Replace-One "445×3=" "297×7="
Replace-One "136×4=" "595×4="
Replace-One "826×7=" "224×2="
Replace-One "862×4=" "188×5="
Replace-One "108×9=" "991×4="
Replace-One "572×8=" "450×9="
Replace-One "569×4=" "200×7="
Replace-One "957×5=" "745×2="
Replace-One "303×6=" "757×2="
Replace-One "628×7=" "714×4="
Replace-One "773×7=" "249×8="
Replace-One "645×3=" "290×5="
Replace-One "275×9=" "739×8="
Replace-One "691×3=" "759×6="
Replace-One "652×2=" "766×6="
